$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '243.60'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '23.83'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.238'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05823'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '6.464'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8084'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.8768'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07271'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.03090'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03054'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09327'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.857'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.04702'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0006020'
$ws.Range('E18').Value = '17OneONEWorstin24h'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.006230'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.004586'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.00008700'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.558'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.183'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1319'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0002340'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03782'
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1053'
$ws.Range('E41').Value = '40BKEXTokenBKK'
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.002570'
$ws.Range('E42').Value = '41CEJICEJI'
$ws.Range('B43').Value = 'KickToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.003241'
$ws.Range('E43').Value = '42KickTokenKICK'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.007777'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005470'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5979'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.02177'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002100'
